$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Transitioning to the new data format: a new project column
# ("Sicangu Online Marketplace") is introduced between the existing
# "ShockTalk" (AE) and "Sisters of Code" (old AF, now shifted to AG)
# columns. Inserting a whole column shifts all the old AF:AQ data
# right by one (to AG:AR) and extends the used range to AR.
$ws.Columns("AF:AF").Insert()

# Header for the newly inserted column.
$ws.Range("AF1").Value = "Sicangu Online Marketplace"

# None of the 18 organizations in this sheet match the new project,
# so the new column is all zeros.
$ws.Range("AF2:AF19").Value = 0

# The row-label header (A1) is renamed from "Org_y" to "Org_x".
$ws.Range("A1").Value = "Org_x"
